$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the used range (rows 298-302) by copying the formatting of the last
# existing row (297) down, then overwrite with the new data values.
$ws.Range("A297:W297").Copy()
$ws.Range("A298:W302").PasteSpecial(-4122)  # xlPasteFormats

# Column B holds dates stored as literal text (e.g. "11/12/2024"); force text
# formatting before assignment so Excel does not auto-convert them to date
# serial numbers.
$ws.Range("B298:B302").NumberFormat = "@"

# Row 298
$ws.Range("A298").Value = 296
$ws.Range("B298").Value = "11/12/2024"
$ws.Range("C298").Value = 138.9583333333333
$ws.Range("D298").Value = 155.2
$ws.Range("E298").Value = 147.57
$ws.Range("F298").Value = 142.5
$ws.Range("G298").Value = 64.75
$ws.Range("H298").Value = 134
$ws.Range("I298").Value = 100.3333333333333
$ws.Range("J298").Value = 100.7964285714286
$ws.Range("K298").Value = 165.6666666666667
$ws.Range("L298").Value = 171.5
$ws.Range("M298").Value = 120
$ws.Range("N298").Value = 172
$ws.Range("O298").Value = 189
$ws.Range("P298").Value = 181
$ws.Range("Q298").Value = 174.5
$ws.Range("R298").Value = 99
$ws.Range("S298").Value = 160
$ws.Range("T298").Value = 0.3192307692307693
$ws.Range("U298").Value = 64.48999999999999
$ws.Range("V298").Value = 118.25
$ws.Range("W298").Value = 64.48999999999999

# Row 299
$ws.Range("A299").Value = 297
$ws.Range("B299").Value = "12/12/2024"
$ws.Range("C299").Value = 132.7222222222222
$ws.Range("D299").Value = 154.2871428571429
$ws.Range("E299").Value = 146.6
$ws.Range("F299").Value = 142
$ws.Range("G299").Value = 64.75
$ws.Range("H299").Value = 134
$ws.Range("I299").Value = 94.49857142857142
$ws.Range("J299").Value = 100.25
$ws.Range("K299").Value = 157.25
$ws.Range("L299").Value = 171.5
$ws.Range("M299").Value = 120
$ws.Range("N299").Value = 172
$ws.Range("O299").Value = 187.8
$ws.Range("P299").Value = 182.7
$ws.Range("Q299").Value = 174.5
$ws.Range("R299").Value = 99
$ws.Range("S299").Value = 160
$ws.Range("T299").Value = 0.3192307692307693
$ws.Range("U299").Value = 64.48999999999999
$ws.Range("V299").Value = 118.25
$ws.Range("W299").Value = 64.48999999999999

# Row 300
$ws.Range("A300").Value = 298
$ws.Range("B300").Value = "13/12/2024"
$ws.Range("C300").Value = 135
$ws.Range("D300").Value = 161.6666666666667
$ws.Range("E300").Value = 147.25
$ws.Range("F300").Value = 142
$ws.Range("G300").Value = 64.75
$ws.Range("H300").Value = 134
$ws.Range("I300").Value = 97.75
$ws.Range("J300").Value = 101.6666666666667
$ws.Range("K300").Value = 164
$ws.Range("L300").Value = 171.5
$ws.Range("M300").Value = 120
$ws.Range("N300").Value = 172
$ws.Range("O300").Value = 187.8
$ws.Range("P300").Value = 182.7
$ws.Range("Q300").Value = 174.5
$ws.Range("R300").Value = 99
$ws.Range("S300").Value = 160
$ws.Range("T300").Value = 0.3192307692307693
$ws.Range("U300").Value = 64.48999999999999
$ws.Range("V300").Value = 118.25
$ws.Range("W300").Value = 64.48999999999999

# Row 301
$ws.Range("A301").Value = 299
$ws.Range("B301").Value = "16/12/2024"
$ws.Range("C301").Value = 130.0625
$ws.Range("D301").Value = 153.5
$ws.Range("E301").Value = 144.5
$ws.Range("F301").Value = 141
$ws.Range("G301").Value = 64.75
$ws.Range("H301").Value = 134
$ws.Range("I301").Value = 88.63047619047619
$ws.Range("J301").Value = 99.366
$ws.Range("K301").Value = 161.5
$ws.Range("L301").Value = 165
$ws.Range("M301").Value = 120
$ws.Range("N301").Value = 172
$ws.Range("O301").Value = 187.8
$ws.Range("P301").Value = 182.7
$ws.Range("Q301").Value = 174.5
$ws.Range("R301").Value = 99
$ws.Range("S301").Value = 160
$ws.Range("T301").Value = 0.3192307692307693
$ws.Range("U301").Value = 64.48999999999999
$ws.Range("V301").Value = 118.25
$ws.Range("W301").Value = 64.48999999999999

# Row 302
$ws.Range("A302").Value = 300
$ws.Range("B302").Value = "17/12/2024"
$ws.Range("C302").Value = 128
$ws.Range("D302").Value = 152.875
$ws.Range("E302").Value = 145.4285714285714
$ws.Range("F302").Value = 141
$ws.Range("G302").Value = 64.75
$ws.Range("H302").Value = 134
$ws.Range("I302").Value = 82.71428571428571
$ws.Range("J302").Value = 97.90909090909091
$ws.Range("K302").Value = 159
$ws.Range("L302").Value = 163.75
$ws.Range("M302").Value = 120
$ws.Range("N302").Value = 162.3333333333333
$ws.Range("O302").Value = 187
$ws.Range("P302").Value = 180
$ws.Range("Q302").Value = 174.5
$ws.Range("R302").Value = 99
$ws.Range("S302").Value = 160
$ws.Range("T302").Value = 0.3192307692307693
$ws.Range("U302").Value = 64.48999999999999
$ws.Range("V302").Value = 118.25
$ws.Range("W302").Value = 64.48999999999999

Write-Host "Added rows 298-302 (A298:W302)."
